# Auto-generated edit script updating cryptos price/volume columns (D, E)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.924.29"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.79%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.843.01"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.08%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.007"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.23%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.19"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.21%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.006"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.12%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4682"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.45%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3631"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.03%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07178"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.44%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9379"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.54%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.63"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.96%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07683"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.60%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.852.11"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.64%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.288"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.10%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.389"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.06%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "88.27"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.45%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.009"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.18%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008585"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.16%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "26.897.44"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.63%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.35"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.01%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.034"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.25%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.64"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.12%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.925"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.85%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.28"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.77%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "18.03"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.14%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.016"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.49%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "114.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.70%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.920"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.22%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08851"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.83%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.165"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.10%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.850"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.32%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.184"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +7.32%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7479"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.86%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.469"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.58%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.087"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.13%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.991"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.35%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01933"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.07%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05152"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.03%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5117"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.05%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.924"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.08%  "

$ws.Range("E42").Value = "  +0.26%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.173"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.05%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4711"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.03%  "

$ws.Range("E45").Value = "  +2.74%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.008"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.22%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "100.03"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.45%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.600"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.66%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06055"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.35%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "64.09"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.57%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "36.17"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.04%  "
